$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 407.66666
$ws.Range("I12").Value = 386
$ws.Range("K12").Value = 386
$ws.Range("M12").Value = -216
$ws.Range("H43").Value = 1950.1818
$ws.Range("I43").Value = 1035
$ws.Range("K43").Value = 1035
$ws.Range("M43").Value = -966
$ws.Range("H62").Value = 4312.6665
$ws.Range("I62").Value = 3169.6
$ws.Range("J62").Value = 6598.8
$ws.Range("K62").Value = 3169.6
$ws.Range("L62").Value = 6598.8
$ws.Range("M62").Value = -2545.6
$ws.Range("N62").Value = -7846.8
$ws.Range("H65").Value = 4312.6665
$ws.Range("I65").Value = 3169.6
$ws.Range("J65").Value = 6598.8
$ws.Range("K65").Value = 15848
$ws.Range("L65").Value = 32994
$ws.Range("M65").Value = -12728
$ws.Range("N65").Value = -39234
$ws.Range("H76").Value = 4300.5
$ws.Range("I76").Value = 3901.2222
$ws.Range("K76").Value = 3901.2222
$ws.Range("M76").Value = -3586.2222
$ws.Range("H79").Value = 4300.5
$ws.Range("I79").Value = 3901.2222
$ws.Range("K79").Value = 3901.2222
$ws.Range("M79").Value = -2809.2222
$ws.Range("H80").Value = 43104010
$ws.Range("I80").Value = 13158143
$ws.Range("J80").Value = 100001150
$ws.Range("K80").Value = 39474429
$ws.Range("L80").Value = 300003450
$ws.Range("M80").Value = -39473431
$ws.Range("N80").Value = -300005446
$ws.Range("H83").Value = 43104010
$ws.Range("I83").Value = 13158143
$ws.Range("J83").Value = 100001150
$ws.Range("K83").Value = 118423287
$ws.Range("L83").Value = 900010350
$ws.Range("M83").Value = -118418295
$ws.Range("N83").Value = -900020334
$ws.Range("H131").Value = 7139.364
$ws.Range("I131").Value = 982.125
$ws.Range("J131").Value = 23558.666
$ws.Range("K131").Value = 2946.375
$ws.Range("L131").Value = 70675.99800000001
$ws.Range("M131").Value = 2093.625
$ws.Range("N131").Value = -80755.99800000001
$ws.Range("H132").Value = 393075.2
$ws.Range("I132").Value = 471774.9
$ws.Range("K132").Value = 1415324.7
$ws.Range("M132").Value = -1412794.7
$ws.Range("H133").Value = 105748.43
$ws.Range("J133").Value = 105748.43
$ws.Range("L133").Value = 105748.43
$ws.Range("N133").Value = -115868.43
$ws.Range("H138").Value = 1729.9899
$ws.Range("J138").Value = 1821.369
$ws.Range("L138").Value = 5464.107
$ws.Range("N138").Value = -15744.107

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2932.47
$ws.Range("I32").Value = 2524.9795
$ws.Range("J32").Value = 22899.5
$ws.Range("K32").Value = 2524.9795
$ws.Range("L32").Value = 22899.5
$ws.Range("M32").Value = -2237.9795
$ws.Range("N32").Value = -23473.5
$ws.Range("H45").Value = 2382.0715
$ws.Range("I45").Value = 2001.4
$ws.Range("J45").Value = 3333.75
$ws.Range("K45").Value = 2001.4
$ws.Range("L45").Value = 3333.75
$ws.Range("M45").Value = -1624.4
$ws.Range("N45").Value = -4087.75
$ws.Range("H61").Value = 7224.1904
$ws.Range("I61").Value = 7142.5264
$ws.Range("K61").Value = 7142.5264
$ws.Range("M61").Value = -6930.5264
$ws.Range("H122").Value = 1971.4857
$ws.Range("I122").Value = 1808.129
$ws.Range("K122").Value = 5424.387
$ws.Range("M122").Value = -2974.387
$ws.Range("H132").Value = 8783.152
$ws.Range("I132").Value = 8435.531000000001
$ws.Range("J132").Value = 10895.615
$ws.Range("K132").Value = 25306.593
$ws.Range("L132").Value = 32686.845
$ws.Range("M132").Value = -22776.593
$ws.Range("N132").Value = -37746.845
$ws.Range("H136").Value = 7224.1904
$ws.Range("I136").Value = 7142.5264
$ws.Range("K136").Value = 21427.5792
$ws.Range("M136").Value = -18877.5792

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3283.7
$ws.Range("I86").Value = 2959.2942
$ws.Range("K86").Value = 2959.2942
$ws.Range("M86").Value = -1836.2942
$ws.Range("H89").Value = 3283.7
$ws.Range("I89").Value = 2959.2942
$ws.Range("K89").Value = 14796.471
$ws.Range("M89").Value = -9180.471
$ws.Range("H94").Value = 1714.7142
$ws.Range("I94").Value = 1932.7693
$ws.Range("K94").Value = 1932.7693
$ws.Range("M94").Value = -1481.7693
$ws.Range("H105").Value = 3724
$ws.Range("I105").Value = 3799
$ws.Range("K105").Value = 3799
$ws.Range("M105").Value = -2052
$ws.Range("H123").Value = 38111.875
$ws.Range("J123").Value = 48979
$ws.Range("L123").Value = 48979
$ws.Range("N123").Value = -58779
$ws.Range("H134").Value = 2988.2354
$ws.Range("I134").Value = 2863.1562
$ws.Range("K134").Value = 8589.4686
$ws.Range("M134").Value = -6054.4686

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20410874
$ws.Range("I31").Value = 26317410
$ws.Range("K31").Value = 26317410
$ws.Range("M31").Value = -26317115
$ws.Range("H34").Value = 20410874
$ws.Range("I34").Value = 26317410
$ws.Range("K34").Value = 26317410
$ws.Range("M34").Value = -26317208
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("N45").ClearContents()
$ws.Range("H48").Value = 30000
$ws.Range("J48").Value = 30000
$ws.Range("L48").Value = 30000
$ws.Range("N48").Value = -30952
$ws.Range("H117").Value = 53496.332
$ws.Range("J117").Value = 53496.332
$ws.Range("L117").Value = 53496.332
$ws.Range("N117").Value = -62674.332
$ws.Range("I134").Value = 2862.7778
$ws.Range("J134").Value = 2372
$ws.Range("K134").Value = 8588.3334
$ws.Range("L134").Value = 7116
$ws.Range("M134").Value = -6053.3334
$ws.Range("N134").Value = -12186

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 45403.383
$ws.Range("J123").Value = 45403.383
$ws.Range("L123").Value = 45403.383
$ws.Range("N123").Value = -50303.383
$ws.Range("H132").Value = 145525.42
$ws.Range("I132").Value = 183396.27
$ws.Range("J132").Value = 6665.6665
$ws.Range("K132").Value = 550188.8099999999
$ws.Range("L132").Value = 19996.9995
$ws.Range("M132").Value = -547658.8099999999
$ws.Range("N132").Value = -25056.9995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1647.9412
$ws.Range("I16").Value = 1145.7037
$ws.Range("J16").Value = 3585.1428
$ws.Range("K16").Value = 1145.7037
$ws.Range("L16").Value = 3585.1428
$ws.Range("M16").Value = -975.7037
$ws.Range("N16").Value = -3925.1428
$ws.Range("H22").Value = 1159
$ws.Range("I22").Value = 971
$ws.Range("K22").Value = 971
$ws.Range("M22").Value = -676
$ws.Range("H26").Value = 44670
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 44670
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -45260
$ws.Range("H27").Value = 1159
$ws.Range("I27").Value = 971
$ws.Range("K27").Value = 971
$ws.Range("M27").Value = -864
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H61").Value = 2383.9565
$ws.Range("I61").Value = 2383.9565
$ws.Range("K61").Value = 2383.9565
$ws.Range("M61").Value = -2181.9565
$ws.Range("H63").Value = 119199.8
$ws.Range("J63").Value = 120000
$ws.Range("L63").Value = 120000
$ws.Range("N63").Value = -121498
$ws.Range("H66").Value = 119199.8
$ws.Range("J66").Value = 120000
$ws.Range("L66").Value = 360000
$ws.Range("N66").Value = -367488
$ws.Range("H74").Value = 79900
$ws.Range("I74").Value = 19750
$ws.Range("K74").Value = 19750
$ws.Range("M74").Value = -18752
$ws.Range("H77").Value = 79900
$ws.Range("I77").Value = 19750
$ws.Range("K77").Value = 59250
$ws.Range("M77").Value = -54258
$ws.Range("H113").Value = 2383.9565
$ws.Range("I113").Value = 2383.9565
$ws.Range("K113").Value = 2383.9565
$ws.Range("M113").Value = -213.9564999999998
$ws.Range("H136").Value = 4279.1665
$ws.Range("I136").Value = 2668.875
$ws.Range("K136").Value = 8006.625
$ws.Range("M136").Value = -5456.625

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 9999
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 9999
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -10579
$ws.Range("H62").Value = 3278.5
$ws.Range("I62").Value = 3278.5
$ws.Range("K62").Value = 3278.5
$ws.Range("M62").Value = -2654.5
$ws.Range("H65").Value = 3278.5
$ws.Range("I65").Value = 3278.5
$ws.Range("K65").Value = 16392.5
$ws.Range("M65").Value = -13272.5
$ws.Range("H81").Value = 12600.52
$ws.Range("I81").Value = 8568.200000000001
$ws.Range("J81").Value = 15288.733
$ws.Range("K81").Value = 17136.4
$ws.Range("L81").Value = 30577.466
$ws.Range("M81").Value = -16075.4
$ws.Range("N81").Value = -32699.466
$ws.Range("H84").Value = 12600.52
$ws.Range("I84").Value = 8568.200000000001
$ws.Range("J84").Value = 15288.733
$ws.Range("K84").Value = 85682
$ws.Range("L84").Value = 152887.33
$ws.Range("M84").Value = -80378
$ws.Range("N84").Value = -163495.33
$ws.Range("H113").Value = 783.6
$ws.Range("I113").Value = 504.15384
$ws.Range("K113").Value = 1512.46152
$ws.Range("M113").Value = 657.5384799999999
$ws.Range("H133").Value = 59170.5
$ws.Range("J133").Value = 59170.5
$ws.Range("L133").Value = 59170.5
$ws.Range("N133").Value = -69290.5
